$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "unna"
$ws.Range("B2").Value = "51.5333,7.6833"

$ws.Range("A2").Select()

$ws.Columns.Item(2).AutoFit()
